$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.177.90"
$ws.Range("E2").Value = "  +0.20%  "

$ws.Range("D3").Value = "1.601.88"
$ws.Range("E3").Value = "  -0.33%  "

$ws.Range("E4").Value = "  +0.24%  "

$ws.Range("D5").Value = "'211.85"
$ws.Range("E5").Value = "  -0.45%  "

$ws.Range("E6").Value = "  +0.20%  "

$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("E8").Value = "  -0.84%  "

$ws.Range("D9").Value = "'0.0614"
$ws.Range("E9").Value = "  -0.88%  "

$ws.Range("D10").Value = "'18.11"
$ws.Range("E10").Value = "  -1.73%  "

$ws.Range("E11").Value = "  +2.41%  "

$ws.Range("D12").Value = "1.822.21"
$ws.Range("E12").Value = "  -0.44%  "

$ws.Range("D13").Value = "1.606.69"
$ws.Range("E13").Value = "  +0.81%  "

$ws.Range("E14").Value = "  -0.85%  "

$ws.Range("D15").Value = "'0.516"
$ws.Range("E15").Value = "  +0.86%  "

$ws.Range("D16").Value = "26.164.99"
$ws.Range("E16").Value = "  +0.18%  "

$ws.Range("D17").Value = "'60.94"
$ws.Range("E17").Value = "  +0.13%  "

$ws.Range("D18").Value = "0.0₃0726"
$ws.Range("E18").Value = "  -0.75%  "

$ws.Range("E19").Value = "  +0.23%  "

$ws.Range("D20").Value = "'204.40"
$ws.Range("E20").Value = "  +3.41%  "

$ws.Range("E21").Value = "  +0.12%  "

$ws.Range("E22").Value = "  -2.42%  "

$ws.Range("E23").Value = "  +0.37%  "

$ws.Range("D24").Value = "'1.94"
$ws.Range("E24").Value = "  +12.67%  "

$ws.Range("D25").Value = "'144.35"
$ws.Range("E25").Value = "  +1.19%  "

$ws.Range("E26").Value = "  +0.31%  "

$ws.Range("E27").Value = "  -7.62%  "

$ws.Range("D28").Value = "'15.21"
$ws.Range("E28").Value = "  -0.09%  "

$ws.Range("D29").Value = "'6.52"
$ws.Range("E29").Value = "  +0.07%  "

$ws.Range("E30").Value = "  +1.72%  "

$ws.Range("E31").Value = "  -0.45%  "

$ws.Range("E32").Value = "  -0.25%  "

$ws.Range("E33").Value = "  -4.38%  "

$ws.Range("D34").Value = "'1.48"
$ws.Range("E34").Value = "  -2.20%  "

$ws.Range("E35").Value = "  +0.19%  "

$ws.Range("D36").Value = "1.140.51"
$ws.Range("E36").Value = "  +3.25%  "

$ws.Range("E37").Value = "  +6.36%  "

$ws.Range("E38").Value = "  +0.26%  "

$ws.Range("E39").Value = "  -1.78%  "

$ws.Range("E40").Value = "  -0.25%  "

$ws.Range("D41").Value = "'0.495"
$ws.Range("E41").Value = "  -2.54%  "

$ws.Range("E42").Value = "  -2.36%  "

$ws.Range("E43").Value = "  +0.40%  "

$ws.Range("D44").Value = "1.737.72"
$ws.Range("E44").Value = "  -0.31%  "

$ws.Range("D45").Value = "'92.13"
$ws.Range("E45").Value = "  -1.12%  "

$ws.Range("E46").Value = "  -2.82%  "

$ws.Range("D47").Value = "'54.07"
$ws.Range("E47").Value = "  +0.38%  "

$ws.Range("E48").Value = "  -0.37%  "

$ws.Range("E49").Value = "  -0.36%  "

$ws.Range("E50").Value = "  +0.43%  "

$ws.Range("D51").Value = "0.0₇0950"
$ws.Range("E51").Value = "  -11.24%  "
